# The commit adds one new weekly price observation for
# "Feria Lagunitas de Puerto Montt - Cebollín" on top of the existing daily
# series: a brand-new row is inserted just above the current row 67,
# pushing every following record down by one row (old row 67 becomes new
# row 68, ..., old row 163 becomes new row 164). The newly inserted row
# duplicates all the fixed attributes of the (old) row 67 record, but
# carries a new date (2021-09-28) and a new Volumen (180) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 67; everything from 67..163 shifts to 68..164.
$ws.Rows.Item(67).Insert()

# Duplicate the record that is now sitting in row 68 (the original row 67)
# into the freshly-inserted row 67, so all the non-changing columns
# (Mercado, Región, Categoría, Calidad, Precio mínimo/máximo/promedio,
# Unidad de comercialización, Origen, Precio $/Kg, Kg o Unidades,
# Clasificación, ...) line up with the rest of the series.
$ws.Rows.Item(68).Copy()
$ws.Rows.Item(67).PasteSpecial()
$excel.CutCopyMode = 0

# Now overwrite the two columns that actually change for this new record:
# Fecha (D, written as the Excel serial date number for 2021-09-28 so the
# existing date number-format on the column renders it correctly) and
# Volumen (J).
$ws.Cells.Item(67, 4).Value = 44467
$ws.Cells.Item(67, 10).Value = 180
